$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Row 1 (header row): F1 becomes a plain number, new columns G1:L1 added
$ws1.Range("F1").Value = 200
$ws1.Range("G1").Value = 115
$ws1.Range("H1").Value = 300
$ws1.Range("I1").Value = 250
$ws1.Range("J1").Value = "250/60"
$ws1.Range("K1").Value = "50/250/60"
$ws1.Range("L1").Value = "50/250/75"

# Row 2
$ws1.Range("G2").Value = 20
$ws1.Range("H2").Value = 2
$ws1.Range("I2").Value = 2
$ws1.Range("J2").Value = 7
$ws1.Range("K2").Value = 16
$ws1.Range("L2").Value = 8

# Row 3
$ws1.Range("G3").Value = 14
$ws1.Range("H3").Value = 2
$ws1.Range("I3").Value = 2
$ws1.Range("J3").Value = 4
$ws1.Range("K3").Value = 7
$ws1.Range("L3").Value = 5

# Row 4
$ws1.Range("G4").Value = 18
$ws1.Range("H4").Value = 6
$ws1.Range("I4").Value = 7
$ws1.Range("J4").Value = 13
$ws1.Range("K4").Value = 22
$ws1.Range("L4").Value = 13

# Row 5
$ws1.Range("G5").Value = 6
$ws1.Range("H5").Value = 2
$ws1.Range("I5").Value = 2
$ws1.Range("J5").Value = 5
$ws1.Range("K5").Value = 8
$ws1.Range("L5").Value = 5

# Row 6
$ws1.Range("G6").Value = 10
$ws1.Range("H6").Value = 0
$ws1.Range("I6").Value = 2
$ws1.Range("J6").Value = 4
$ws1.Range("K6").Value = 10
$ws1.Range("L6").Value = 7

# Row 7
$ws1.Range("G7").Value = 10
$ws1.Range("H7").Value = 2
$ws1.Range("I7").Value = 4
$ws1.Range("J7").Value = 6
$ws1.Range("K7").Value = 7
$ws1.Range("L7").Value = 5

# Row 8
$ws1.Range("G8").Value = 8
$ws1.Range("H8").Value = 0
$ws1.Range("I8").Value = 1
$ws1.Range("J8").Value = 3
$ws1.Range("K8").Value = 6
$ws1.Range("L8").Value = 3

# Row 9
$ws1.Range("G9").Value = 3
$ws1.Range("H9").Value = 0
$ws1.Range("I9").Value = 0
$ws1.Range("J9").Value = 0
$ws1.Range("K9").Value = 1
$ws1.Range("L9").Value = 0

# Row 10
$ws1.Range("G10").Value = 0
$ws1.Range("H10").Value = 0
$ws1.Range("I10").Value = 0
$ws1.Range("J10").Value = 0
$ws1.Range("K10").Value = 0
$ws1.Range("L10").Value = 0

# Row 11
$ws1.Range("G11").Value = 5
$ws1.Range("H11").Value = 0
$ws1.Range("I11").Value = 0
$ws1.Range("J11").Value = 0
$ws1.Range("K11").Value = 0
$ws1.Range("L11").Value = 0

# Make Sheet1 the active / selected tab, matching the author's workflow
[void]$ws1.Select()
[void]$ws1.Range("L12").Select()
